$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update J1 anchor-score label cell (word stays "positive", only its
#     underlying shared-string slot changes upstream) ---
$ws.Range("J1").Value = "positive"

# --- Update A column word labels for rows 3-6 (dataset grew, the anchor
#     words "sc" and "panic" swap order) ---
$ws.Range("A3").Value = "crude"
$ws.Range("A4").Value = "crisis"
$ws.Range("A5").Value = "sc"
$ws.Range("A6").Value = "panic"

# --- Update numeric stats columns B:H for rows 3-6 (anchor-word table) ---
$ws.Range("B3").Value = 0.7647058823529411
$ws.Range("C3").Value = 26
$ws.Range("D3").Value = 26
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 8
$ws.Range("B4").Value = 0.6095890410958904
$ws.Range("C4").Value = 178
$ws.Range("D4").Value = 178
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 114
$ws.Range("B5").Value = 0.2063492063492063
$ws.Range("C5").Value = 39
$ws.Range("D5").Value = 39
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 150
$ws.Range("B6").Value = 0.1744186046511628
$ws.Range("C6").Value = 90
$ws.Range("D6").Value = 90
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 426

# --- Update word labels for the confidence-ranked table (J column), rows 3-27 ---
#     These keep their row position but the word each one shows changes because
#     the underlying word list order changed upstream.
$ws.Range("J3").Value = "happy"
$ws.Range("J4").Value = "interesting"
$ws.Range("J5").Value = "best"
$ws.Range("J6").Value = "love"
$ws.Range("J7").Value = "great"
$ws.Range("J8").Value = "positive"
$ws.Range("J9").Value = "won"
$ws.Range("J10").Value = "thank"
$ws.Range("J11").Value = "thanks"
$ws.Range("J12").Value = "special"
$ws.Range("J13").Value = "free"
$ws.Range("J14").Value = "safe"
$ws.Range("J15").Value = "good"
$ws.Range("J16").Value = "heroes"
$ws.Range("J17").Value = "support"
$ws.Range("J18").Value = "safety"
$ws.Range("J19").Value = "well"
$ws.Range("J20").Value = "fresh"
$ws.Range("J21").Value = "better"
$ws.Range("J22").Value = "relief"
$ws.Range("J23").Value = "hand"
$ws.Range("J24").Value = "like"
$ws.Range("J25").Value = "care"
$ws.Range("J26").Value = "help"
$ws.Range("J27").Value = "protect"

# --- Update numeric stats columns K:Q for rows 3-27 ---
$ws.Range("K3").Value = 0.9615384615384616
$ws.Range("L3").Value = 25
$ws.Range("M3").Value = 25
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 1
$ws.Range("K4").Value = 0.9393939393939394
$ws.Range("L4").Value = 31
$ws.Range("M4").Value = 31
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 2
$ws.Range("K5").Value = 0.9152542372881356
$ws.Range("L5").Value = 54
$ws.Range("M5").Value = 54
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 5
$ws.Range("K6").Value = 0.8913043478260869
$ws.Range("L6").Value = 41
$ws.Range("M6").Value = 41
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 5
$ws.Range("K7").Value = 0.8839285714285714
$ws.Range("L7").Value = 99
$ws.Range("M7").Value = 99
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 13
$ws.Range("K8").Value = 0.8275862068965517
$ws.Range("L8").Value = 48
$ws.Range("M8").Value = 48
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 10
$ws.Range("K9").Value = 0.8205128205128205
$ws.Range("L9").Value = 32
$ws.Range("M9").Value = 32
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 7
$ws.Range("K10").Value = 0.78125
$ws.Range("L10").Value = 100
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 28
$ws.Range("K11").Value = 0.7804878048780488
$ws.Range("L11").Value = 64
$ws.Range("M11").Value = 64
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 18
$ws.Range("K12").Value = 0.7777777777777778
$ws.Range("L12").Value = 28
$ws.Range("M12").Value = 28
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 8
$ws.Range("K13").Value = 0.7666666666666667
$ws.Range("L13").Value = 92
$ws.Range("M13").Value = 92
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 28
$ws.Range("K14").Value = 0.7323943661971831
$ws.Range("L14").Value = 104
$ws.Range("M14").Value = 104
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 38
$ws.Range("K15").Value = 0.6875
$ws.Range("L15").Value = 110
$ws.Range("M15").Value = 110
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 50
$ws.Range("K16").Value = 0.6808510638297872
$ws.Range("L16").Value = 32
$ws.Range("M16").Value = 32
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 15
$ws.Range("K17").Value = 0.6698113207547169
$ws.Range("L17").Value = 71
$ws.Range("M17").Value = 71
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 35
$ws.Range("K18").Value = 0.6274509803921569
$ws.Range("L18").Value = 32
$ws.Range("M18").Value = 32
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 19
$ws.Range("K19").Value = 0.6170212765957447
$ws.Range("L19").Value = 58
$ws.Range("M19").Value = 58
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 36
$ws.Range("K20").Value = 0.5833333333333334
$ws.Range("L20").Value = 28
$ws.Range("M20").Value = 28
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 20
$ws.Range("K21").Value = 0.5714285714285714
$ws.Range("L21").Value = 36
$ws.Range("M21").Value = 36
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 27
$ws.Range("K22").Value = 0.56
$ws.Range("L22").Value = 28
$ws.Range("M22").Value = 28
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 22
$ws.Range("K23").Value = 0.4934725848563969
$ws.Range("L23").Value = 189
$ws.Range("M23").Value = 189
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 194
$ws.Range("K24").Value = 0.4470588235294118
$ws.Range("L24").Value = 152
$ws.Range("M24").Value = 152
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 188
$ws.Range("K25").Value = 0.4269662921348314
$ws.Range("L25").Value = 38
$ws.Range("M25").Value = 38
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 51
$ws.Range("K26").Value = 0.4101694915254237
$ws.Range("L26").Value = 121
$ws.Range("M26").Value = 121
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 174
$ws.Range("K27").Value = 0.3698630136986301
$ws.Range("L27").Value = 27
$ws.Range("M27").Value = 27
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 46

# --- Append new rows 28-35 (larger dataset produced more ranked words) ---
# Copy the formatting (bold + border + alignment) from row 27 down through row 35
# so the new rows visually match the existing ranked-word rows.
$ws.Range("J27:Q27").Copy()
$ws.Range("J28:Q35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("J28").Value = "increase"
$ws.Range("K28").Value = 0.3461538461538461
$ws.Range("L28").Value = 27
$ws.Range("M28").Value = 27
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 51
$ws.Range("J29").Value = "please"
$ws.Range("K29").Value = 0.2845188284518829
$ws.Range("L29").Value = 68
$ws.Range("M29").Value = 68
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 171
$ws.Range("J30").Value = "you"
$ws.Range("K30").Value = 0.0275
$ws.Range("L30").Value = 33
$ws.Range("M30").Value = 33
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = $false
$ws.Range("Q30").Value = 1167
$ws.Range("J31").Value = "and"
$ws.Range("K31").Value = 0.01161048689138577
$ws.Range("L31").Value = 31
$ws.Range("M31").Value = 34
$ws.Range("N31").Value = 0.91
$ws.Range("O31").Value = 0.08999999999999997
$ws.Range("P31").Value = $true
$ws.Range("Q31").Value = 2639
$ws.Range("J32").Value = ","
$ws.Range("K32").Value = 0.01145662847790507
$ws.Range("L32").Value = 28
$ws.Range("M32").Value = 29
$ws.Range("N32").Value = 0.97
$ws.Range("O32").Value = 0.03000000000000003
$ws.Range("P32").Value = $true
$ws.Range("Q32").Value = 2416
$ws.Range("J33").Value = "to"
$ws.Range("K33").Value = 0.006700554528650647
$ws.Range("L33").Value = 29
$ws.Range("M33").Value = 29
$ws.Range("N33").Value = 1
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = $false
$ws.Range("Q33").Value = 4299
$ws.Range("J34").Value = "the"
$ws.Range("K34").Value = 0.006005424254165052
$ws.Range("L34").Value = 31
$ws.Range("M34").Value = 34
$ws.Range("N34").Value = 0.91
$ws.Range("O34").Value = 0.08999999999999997
$ws.Range("P34").Value = $true
$ws.Range("Q34").Value = 5131
$ws.Range("J35").Value = "."
$ws.Range("K35").Value = 0.005604483586869495
$ws.Range("L35").Value = 28
$ws.Range("M35").Value = 28
$ws.Range("N35").Value = 1
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = $false
$ws.Range("Q35").Value = 4968

